$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 41 : Finalisation des endpoints de sélection d'utilisateurs ---
$ws.Range("A41").Value = "Réalisation "
$ws.Rows.Item(41).RowHeight = 75

# --- Fix existing row 40 (Column B) text: "un utilisateur" -> "des utilisateurs" ---
$ws.Range("B40").Value = "Création du endpoint API qui permet de sélectionner des utilisateurs"

$ws.Range("B41").Value = "Finalisation du endpoint de selection des utilisateurs et création du endpoints permettant de sélectionner un seul utilisateur"
$ws.Range("C41").Value = 2
$ws.Range("D41").Value = (Get-Date -Year 2019 -Month 2 -Day 27).Date

# --- Row 42 : Endpoint de gestion des tokens ---
$ws.Range("A42").Value = "Réalisation"
$ws.Range("A42").WrapText = $true
$ws.Range("B42").Value = "Création du endpoint de gestion des tokens et réalisation de toute la logique de création de tokens"
$ws.Range("C42").Value = 2.5
$ws.Range("D42").Value = (Get-Date -Year 2019 -Month 3 -Day 4).Date
$ws.Rows.Item(42).RowHeight = 60

# --- Row 43 : Tests sur le endpoint de récupération de tokens ---
$ws.Range("A43").Value = "Réalisation"
$ws.Range("A43").WrapText = $true
$ws.Range("B43").Value = "Tests sur le endpoint de récupération de tokens pour trouver des bugs. Un bug était présent lorsque aucun paramètre n'était envoyé avec la requête"
$ws.Range("C43").Value = 0.75
$ws.Range("D43").Value = (Get-Date -Year 2019 -Month 3 -Day 5).Date
$ws.Rows.Item(43).RowHeight = 90

# --- Row 44 : Mise à jour du script SQL pour compatibilité MySQL ---
$ws.Range("A44").Value = "Réalisation"
$ws.Range("A44").WrapText = $true
$ws.Range("B44").Value = "Mise à jour du script SQL de création de la base de données pour rendre le script compatible avec les anciennes versions de MYSQL"
$ws.Range("C44").Value = 0.5
$ws.Range("D44").Value = (Get-Date -Year 2019 -Month 3 -Day 6).Date
$ws.Rows.Item(44).RowHeight = 75

# --- Row 45 : Refactoring de la gestion des JWT ---
$ws.Range("A45").Value = "Réalisation"
$ws.Range("A45").WrapText = $true
$ws.Range("B45").Value = "Refactoring de la gestion des JWT pour améliorer la lisibilité et la réutilisation du code. Création de la fonction de vérification des JWT pour authentifier les utilisateurs"
$ws.Range("C45").Value = 1
$ws.Range("D45").Value = (Get-Date -Year 2019 -Month 3 -Day 6).Date
$ws.Rows.Item(45).RowHeight = 90

# --- Update frozen pane / selection to reflect newly scrolled view ---
$ws.Range("A44").Select()
$excel.ActiveWindow.ScrollRow = 44
$ws.Range("C46").Select()
